$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21)
$data = @(
    @(-16.96957763879266, -0.4378763594642169, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266),
    @(-16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266),
    @(-16.96957763879266, -0.4852424790317624, 0.357808341930667, -16.96957763879266, 3.936741219072015, -16.96957763879266, 0.04366869261774536, -16.96957763879266, 0.9466808923777491, -16.96957763879266),
    @(-16.96957763879266, 0.17572764847306, -16.96957763879266, -16.96957763879266, -16.96957763879266, 3.735836228783353, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266),
    @(-16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266),
    @(3.036525138304221, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266),
    @(-16.96957763879266, -16.96957763879266, -16.96957763879266, 2.266737924218202, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266),
    @(3.560071189740182, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266),
    @(-16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, 1.63905620474978),
    @(-16.96957763879266, -16.96957763879266, -16.96957763879266, 2.011624486069923, -16.96957763879266, 1.569251638115098, -16.96957763879266, -16.96957763879266, -16.96957763879266, 0.8581682114067475),
    @(-16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266),
    @(-16.96957763879266, -16.96957763879266, -16.96957763879266, 1.696511105010069, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, 1.776323138838111, 1.048888803633094),
    @(-16.96957763879266, -16.96957763879266, 1.984098768686467, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, 1.550645848544933),
    @(-16.96957763879266, -16.96957763879266, -0.3155990051630655, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266),
    @(-16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, 3.25029160656327, -16.96957763879266),
    @(-16.96957763879266, 0.4251613628003216, -0.08301177380481602, -16.96957763879266, -16.96957763879266, -16.96957763879266, 3.091235694940845, -16.96957763879266, 1.655492607626888, -16.96957763879266),
    @(-16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, -16.96957763879266, 0.95215290743734, -16.96957763879266, 0.9862595399720255, -16.96957763879266),
    @(-16.96957763879266, -16.96957763879266, 2.795004049420683, -16.96957763879266, -16.96957763879266, -16.96957763879266, 1.547649604062878, 4.321917415462399, -16.96957763879266, -16.96957763879266),
    @(-16.96957763879266, 3.032244226477448, 2.602724752671468, -16.96957763879266, 2.228438079457155, -16.96957763879266, 1.721542225930646, -16.96957763879266, -16.96957763879266, 3.332580422531509),
    @(-16.96957763879266, 2.980797420243182, -16.96957763879266, 2.984409880644038, -16.96957763879266, 1.891209280519887, 1.195883897823433, -16.96957763879266, -16.96957763879266, -16.96957763879266)
)

for ($ri = 0; $ri -lt $rows.Length; $ri++) {
    $r = $rows[$ri]
    $vals = $data[$ri]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
    }
}
